$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------
# /bin row (row 2): clarify that /bin holds executables used by ALL users
$ws.Range("B2").Value = "исполняемые файлы, используемые всеми пользователями"

# /sbin row (row 3): clarify that /sbin holds executables for configuring the OS
$ws.Range("B3").Value = "исполняемые файлы для настройки ОС"

# /opt row (row 14): fix typo "под-каталоги" -> "подкаталоги"
$ws.Range("D14").Value = "в эту папку устанавливаются проприетарные программы, игры или драйвера.`nЭто программы созданные в виде отдельных исполняемых файлов самими производителями. `nТакие программы устанавливаются в подкаталоги /opt/ `nвсе исполняемые файлы, библиотеки и файлы конфигурации находятся в одной папке."

# --- View / selection changes --------------------------------------------
$window = $excel.ActiveWindow
$window.View = $excel.Constants.xlNormalView
$window.Zoom = 145
$window.ScrollRow = 13
$window.ScrollColumn = 1
$ws.Range("D15").Select()
